# issue #5: stock data output to json file
#
# Adds a "property_category" column (value "stock") to the 股票 (stock)
# sheet so the JSON exporter can tag each row's property category, and
# cleans up a few shared-string typos (stray full-width spaces inside
# company names, a full-width comma in a share-count value).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("股票")

# Insert a new column at H (before the existing "date" column) and copy
# the header/data cell formatting from the column that is being pushed
# right (old H, now I) so the new column matches the sheet's look.
$ws.Columns("H").Insert()

$ws.Range("I1").Copy() | Out-Null
$ws.Range("H1").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(1, 8).Value = "property_category"

for ($r = 2; $r -le 12; $r++) {
    $ws.Range("I$r").Copy() | Out-Null
    $ws.Range("H$r").PasteSpecial(-4122) | Out-Null
    $ws.Cells.Item($r, 8).Value = "stock"
}

$ws.Application.CutCopyMode = $false

# Fix shared-string typos (stray full-width space / full-width comma).
$ws.Range("B3").Value = "國泰金融控股股份有限公司"
$ws.Range("B5").Value = "台新金融控股股份有限公司"
$ws.Range("D5").Value = "6526"
$ws.Range("B9").Value = "台灣積體電路製造股份有限公司"
$ws.Range("B10").Value = "兆豐票券金融股份有限公司"
$ws.Range("B12").Value = "台灣積體電路製造股份有限公司"
